{"js": "// Update the date paragraph (first paragraph in the body, above the table).\nconst dateResults = context.document.body.paragraphs.getFirst().search(\"2025-02-20 Thursday\", { matchCase: true });\ndateResults.load(\"text\");\nawait context.sync();\nif (dateResults.items.length !== 1) {\n  throw new Error(\"Expected exactly 1 match for date text, found \" + dateResults.items.length);\n}\ndateResults.items[0].insertText(\"2025-02-21 Friday\", \"Replace\");\nawait context.sync();\n\n// Update each arithmetic-problem cell in the table (20 rows x 5 columns).\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\nconst cellUpdates = [\n  { row: 0, col: 0, oldText: \"96-53=\", newText: \"60+13=\" },\n  { row: 0, col: 1, oldText: \"13+68=\", newText: \"7+65=\" },\n  { row: 0, col: 2, oldText: \"60+28=\", newText: \"16+77=\" },\n  { row: 0, col: 3, oldText: \"0+99=\", newText: \"96-90=\" },\n  { row: 0, col: 4, oldText: \"15+8=\", newText: \"35+29=\" },\n  { row: 1, col: 0, oldText: \"9+84=\", newText: \"2+16=\" },\n  { row: 1, col: 1, oldText: \"87-21=\", newText: \"27+15=\" },\n  { row: 1, col: 2, oldText: \"0+22=\", newText: \"45-17=\" },\n  { row: 1, col: 3, oldText: \"69-9=\", newText: \"44-38=\" },\n  { row: 1, col: 4, oldText: \"44-23=\", newText: \"61-56=\" },\n  { row: 2, col: 0, oldText: \"20-13=\", newText: \"5+39=\" },\n  { row: 2, col: 1, oldText: \"54+31=\", newText: \"22+57=\" },\n  { row: 2, col: 2, oldText: \"43-35=\", newText: \"37+55=\" },\n  { row: 2, col: 3, oldText: \"23+0=\", newText: \"65+0=\" },\n  { row: 2, col: 4, oldText: \"11+39=\", newText: \"75-72=\" },\n  { row: 3, col: 0, oldText: \"66-42=\", newText: \"10+51=\" },\n  { row: 3, col: 1, oldText: \"85-19=\", newText: \"14+38=\" },\n  { row: 3, col: 2, oldText: \"89-71=\", newText: \"41+6=\" },\n  { row: 3, col: 3, oldText: \"35+64=\", newText: \"82+7=\" },\n  { row: 3, col: 4, oldText: \"53-24=\", newText: \"79-0=\" },\n  { row: 4, col: 0, oldText: \"87-52=\", newText: \"45+12=\" },\n  { row: 4, col: 1, oldText: \"71-15=\", newText: \"54+35=\" },\n  { row: 4, col: 2, oldText: \"71-63=\", newText: \"32+2=\" },\n  { row: 4, col: 3, oldText: \"22+75=\", newText: \"77-53=\" },\n  { row: 4, col: 4, oldText: \"18+42=\", newText: \"58+27=\" },\n  { row: 5, col: 0, oldText: \"56+13=\", newText: \"21+4=\" },\n  { row: 5, col: 1, oldText: \"72-41=\", newText: \"49+10=\" },\n  { row: 5, col: 2, oldText: \"74+17=\", newText: \"17+16=\" },\n  { row: 5, col: 3, oldText: \"75-7=\", newText: \"67-58=\" },\n  { row: 5, col: 4, oldText: \"6+61=\", newText: \"46+30=\" },\n  { row: 6, col: 0, oldText: \"88-63=\", newText: \"70-33=\" },\n  { row: 6, col: 1, oldText: \"55+28=\", newText: \"3+76=\" },\n  { row: 6, col: 2, oldText: \"29+22=\", newText: \"64+32=\" },\n  { row: 6, col: 3, oldText: \"35-17=\", newText: \"72-20=\" },\n  { row: 6, col: 4, oldText: \"85-75=\", newText: \"92-50=\" },\n  { row: 7, col: 0, oldText: \"14+22=\", newText: \"2+79=\" },\n  { row: 7, col: 1, oldText: \"52+37=\", newText: \"5+2=\" },\n  { row: 7, col: 2, oldText: \"85+4=\", newText: \"22-6=\" },\n  { row: 7, col: 3, oldText: \"56+1=\", newText: \"47-22=\" },\n  { row: 7, col: 4, oldText: \"61-41=\", newText: \"88-86=\" },\n  { row: 8, col: 0, oldText: \"97-42=\", newText: \"22+68=\" },\n  { row: 8, col: 1, oldText: \"67+7=\", newText: \"13+26=\" },\n  { row: 8, col: 2, oldText: \"5-2=\", newText: \"95-60=\" },\n  { row: 8, col: 3, oldText: \"94-81=\", newText: \"7+70=\" },\n  { row: 8, col: 4, oldText: \"14+28=\", newText: \"93-88=\" },\n  { row: 9, col: 0, oldText: \"61+16=\", newText: \"4+77=\" },\n  { row: 9, col: 1, oldText: \"37-26=\", newText: \"94-36=\" },\n  { row: 9, col: 2, oldText: \"61+30=\", newText: \"52+38=\" },\n  { row: 9, col: 3, oldText: \"35-12=\", newText: \"75-29=\" },\n  { row: 9, col: 4, oldText: \"98-20=\", newText: \"93-28=\" },\n  { row: 10, col: 0, oldText: \"23+6=\", newText: \"29+5=\" },\n  { row: 10, col: 1, oldText: \"79-7=\", newText: \"18+70=\" },\n  { row: 10, col: 2, oldText: \"43+20=\", newText: \"52-48=\" },\n  { row: 10, col: 3, oldText: \"47+6=\", newText: \"12+19=\" },\n  { row: 10, col: 4, oldText: \"32-19=\", newText: \"9+35=\" },\n  { row: 11, col: 0, oldText: \"33+40=\", newText: \"44-3=\" },\n  { row: 11, col: 1, oldText: \"21+74=\", newText: \"95-7=\" },\n  { row: 11, col: 2, oldText: \"81-71=\", newText: \"65-61=\" },\n  { row: 11, col: 3, oldText: \"42+15=\", newText: \"85+2=\" },\n  { row: 11, col: 4, oldText: \"10+5=\", newText: \"98-98=\" },\n  { row: 12, col: 0, oldText: \"62+5=\", newText: \"89-41=\" },\n  { row: 12, col: 1, oldText: \"54+2=\", newText: \"30+2=\" },\n  { row: 12, col: 2, oldText: \"38+21=\", newText: \"19+37=\" },\n  { row: 12, col: 3, oldText: \"22+73=\", newText: \"99-44=\" },\n  { row: 12, col: 4, oldText: \"18-5=\", newText: \"4+94=\" },\n  { row: 13, col: 0, oldText: \"35-15=\", newText: \"3+58=\" },\n  { row: 13, col: 1, oldText: \"75-69=\", newText: \"10+55=\" },\n  { row: 13, col: 2, oldText: \"38-10=\", newText: \"77-64=\" },\n  { row: 13, col: 3, oldText: \"36+41=\", newText: \"94-20=\" },\n  { row: 13, col: 4, oldText: \"58+34=\", newText: \"45+30=\" },\n  { row: 14, col: 0, oldText: \"75+3=\", newText: \"16+7=\" },\n  { row: 14, col: 1, oldText: \"52+42=\", newText: \"99-38=\" },\n  { row: 14, col: 2, oldText: \"30+22=\", newText: \"96-41=\" },\n  { row: 14, col: 3, oldText: \"85-39=\", newText: \"40+13=\" },\n  { row: 14, col: 4, oldText: \"61-59=\", newText: \"55+13=\" },\n  { row: 15, col: 0, oldText: \"34+8=\", newText: \"32+51=\" },\n  { row: 15, col: 1, oldText: \"77+21=\", newText: \"73-59=\" },\n  { row: 15, col: 2, oldText: \"89-77=\", newText: \"88-30=\" },\n  { row: 15, col: 3, oldText: \"32+32=\", newText: \"26+13=\" },\n  { row: 15, col: 4, oldText: \"61+17=\", newText: \"51-18=\" },\n  { row: 16, col: 0, oldText: \"2+48=\", newText: \"36+40=\" },\n  { row: 16, col: 1, oldText: \"38-25=\", newText: \"11+1=\" },\n  { row: 16, col: 2, oldText: \"68-62=\", newText: \"36+62=\" },\n  { row: 16, col: 3, oldText: \"71-12=\", newText: \"23-18=\" },\n  { row: 16, col: 4, oldText: \"98-5=\", newText: \"98-93=\" },\n  { row: 17, col: 0, oldText: \"59-50=\", newText: \"14+19=\" },\n  { row: 17, col: 1, oldText: \"83-82=\", newText: \"11+8=\" },\n  { row: 17, col: 2, oldText: \"55+12=\", newText: \"10+3=\" },\n  { row: 17, col: 3, oldText: \"97-51=\", newText: \"96-52=\" },\n  { row: 17, col: 4, oldText: \"58+19=\", newText: \"20-13=\" },\n  { row: 18, col: 0, oldText: \"86-7=\", newText: \"55-0=\" },\n  { row: 18, col: 1, oldText: \"65+18=\", newText: \"22-15=\" },\n  { row: 18, col: 2, oldText: \"71-64=\", newText: \"51-29=\" },\n  { row: 18, col: 3, oldText: \"64-58=\", newText: \"31+30=\" },\n  { row: 18, col: 4, oldText: \"82-33=\", newText: \"63-4=\" },\n  { row: 19, col: 0, oldText: \"98-59=\", newText: \"23+52=\" },\n  { row: 19, col: 1, oldText: \"1+28=\", newText: \"28-7=\" },\n  { row: 19, col: 2, oldText: \"42+38=\", newText: \"15-14=\" },\n  { row: 19, col: 3, oldText: \"59+39=\", newText: \"42-39=\" },\n  { row: 19, col: 4, oldText: \"55-50=\", newText: \"92-85=\" },\n];\n\nfor (const { row, col, oldText, newText } of cellUpdates) {\n  const cell = table.getCell(row, col);\n  const results = cell.body.search(oldText, { matchCase: true });\n  results.load(\"text\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\n      \"Expected exactly 1 match in cell (\" + row + \",\" + col + \") for '\" + oldText + \"', found \" + results.items.length\n    );\n  }\n  results.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Update the date paragraph (first paragraph in the body, above the table).\n$dateParagraph = $d.Paragraphs.Item(1)\nif ($dateParagraph.Range.Text.TrimEnd([char]13,[char]7) -ne \"2025-02-20 Thursday\") {\n    throw \"Unexpected date paragraph text: \" + $dateParagraph.Range.Text\n}\n$dateParagraph.Range.Text = \"2025-02-21 Friday\"\n\n# Update each arithmetic-problem cell in the table (20 rows x 5 columns).\n$table = $d.Tables.Item(1)\n\n$cellUpdates = @(\n    @{ Row = 1; Col = 1; OldText = \"96-53=\"; NewText = \"60+13=\" },\n    @{ Row = 1; Col = 2; OldText = \"13+68=\"; NewText = \"7+65=\" },\n    @{ Row = 1; Col = 3; OldText = \"60+28=\"; NewText = \"16+77=\" },\n    @{ Row = 1; Col = 4; OldText = \"0+99=\"; NewText = \"96-90=\" },\n    @{ Row = 1; Col = 5; OldText = \"15+8=\"; NewText = \"35+29=\" },\n    @{ Row = 2; Col = 1; OldText = \"9+84=\"; NewText = \"2+16=\" },\n    @{ Row = 2; Col = 2; OldText = \"87-21=\"; NewText = \"27+15=\" },\n    @{ Row = 2; Col = 3; OldText = \"0+22=\"; NewText = \"45-17=\" },\n    @{ Row = 2; Col = 4; OldText = \"69-9=\"; NewText = \"44-38=\" },\n    @{ Row = 2; Col = 5; OldText = \"44-23=\"; NewText = \"61-56=\" },\n    @{ Row = 3; Col = 1; OldText = \"20-13=\"; NewText = \"5+39=\" },\n    @{ Row = 3; Col = 2; OldText = \"54+31=\"; NewText = \"22+57=\" },\n    @{ Row = 3; Col = 3; OldText = \"43-35=\"; NewText = \"37+55=\" },\n    @{ Row = 3; Col = 4; OldText = \"23+0=\"; NewText = \"65+0=\" },\n    @{ Row = 3; Col = 5; OldText = \"11+39=\"; NewText = \"75-72=\" },\n    @{ Row = 4; Col = 1; OldText = \"66-42=\"; NewText = \"10+51=\" },\n    @{ Row = 4; Col = 2; OldText = \"85-19=\"; NewText = \"14+38=\" },\n    @{ Row = 4; Col = 3; OldText = \"89-71=\"; NewText = \"41+6=\" },\n    @{ Row = 4; Col = 4; OldText = \"35+64=\"; NewText = \"82+7=\" },\n    @{ Row = 4; Col = 5; OldText = \"53-24=\"; NewText = \"79-0=\" },\n    @{ Row = 5; Col = 1; OldText = \"87-52=\"; NewText = \"45+12=\" },\n    @{ Row = 5; Col = 2; OldText = \"71-15=\"; NewText = \"54+35=\" },\n    @{ Row = 5; Col = 3; OldText = \"71-63=\"; NewText = \"32+2=\" },\n    @{ Row = 5; Col = 4; OldText = \"22+75=\"; NewText = \"77-53=\" },\n    @{ Row = 5; Col = 5; OldText = \"18+42=\"; NewText = \"58+27=\" },\n    @{ Row = 6; Col = 1; OldText = \"56+13=\"; NewText = \"21+4=\" },\n    @{ Row = 6; Col = 2; OldText = \"72-41=\"; NewText = \"49+10=\" },\n    @{ Row = 6; Col = 3; OldText = \"74+17=\"; NewText = \"17+16=\" },\n    @{ Row = 6; Col = 4; OldText = \"75-7=\"; NewText = \"67-58=\" },\n    @{ Row = 6; Col = 5; OldText = \"6+61=\"; NewText = \"46+30=\" },\n    @{ Row = 7; Col = 1; OldText = \"88-63=\"; NewText = \"70-33=\" },\n    @{ Row = 7; Col = 2; OldText = \"55+28=\"; NewText = \"3+76=\" },\n    @{ Row = 7; Col = 3; OldText = \"29+22=\"; NewText = \"64+32=\" },\n    @{ Row = 7; Col = 4; OldText = \"35-17=\"; NewText = \"72-20=\" },\n    @{ Row = 7; Col = 5; OldText = \"85-75=\"; NewText = \"92-50=\" },\n    @{ Row = 8; Col = 1; OldText = \"14+22=\"; NewText = \"2+79=\" },\n    @{ Row = 8; Col = 2; OldText = \"52+37=\"; NewText = \"5+2=\" },\n    @{ Row = 8; Col = 3; OldText = \"85+4=\"; NewText = \"22-6=\" },\n    @{ Row = 8; Col = 4; OldText = \"56+1=\"; NewText = \"47-22=\" },\n    @{ Row = 8; Col = 5; OldText = \"61-41=\"; NewText = \"88-86=\" },\n    @{ Row = 9; Col = 1; OldText = \"97-42=\"; NewText = \"22+68=\" },\n    @{ Row = 9; Col = 2; OldText = \"67+7=\"; NewText = \"13+26=\" },\n    @{ Row = 9; Col = 3; OldText = \"5-2=\"; NewText = \"95-60=\" },\n    @{ Row = 9; Col = 4; OldText = \"94-81=\"; NewText = \"7+70=\" },\n    @{ Row = 9; Col = 5; OldText = \"14+28=\"; NewText = \"93-88=\" },\n    @{ Row = 10; Col = 1; OldText = \"61+16=\"; NewText = \"4+77=\" },\n    @{ Row = 10; Col = 2; OldText = \"37-26=\"; NewText = \"94-36=\" },\n    @{ Row = 10; Col = 3; OldText = \"61+30=\"; NewText = \"52+38=\" },\n    @{ Row = 10; Col = 4; OldText = \"35-12=\"; NewText = \"75-29=\" },\n    @{ Row = 10; Col = 5; OldText = \"98-20=\"; NewText = \"93-28=\" },\n    @{ Row = 11; Col = 1; OldText = \"23+6=\"; NewText = \"29+5=\" },\n    @{ Row = 11; Col = 2; OldText = \"79-7=\"; NewText = \"18+70=\" },\n    @{ Row = 11; Col = 3; OldText = \"43+20=\"; NewText = \"52-48=\" },\n    @{ Row = 11; Col = 4; OldText = \"47+6=\"; NewText = \"12+19=\" },\n    @{ Row = 11; Col = 5; OldText = \"32-19=\"; NewText = \"9+35=\" },\n    @{ Row = 12; Col = 1; OldText = \"33+40=\"; NewText = \"44-3=\" },\n    @{ Row = 12; Col = 2; OldText = \"21+74=\"; NewText = \"95-7=\" },\n    @{ Row = 12; Col = 3; OldText = \"81-71=\"; NewText = \"65-61=\" },\n    @{ Row = 12; Col = 4; OldText = \"42+15=\"; NewText = \"85+2=\" },\n    @{ Row = 12; Col = 5; OldText = \"10+5=\"; NewText = \"98-98=\" },\n    @{ Row = 13; Col = 1; OldText = \"62+5=\"; NewText = \"89-41=\" },\n    @{ Row = 13; Col = 2; OldText = \"54+2=\"; NewText = \"30+2=\" },\n    @{ Row = 13; Col = 3; OldText = \"38+21=\"; NewText = \"19+37=\" },\n    @{ Row = 13; Col = 4; OldText = \"22+73=\"; NewText = \"99-44=\" },\n    @{ Row = 13; Col = 5; OldText = \"18-5=\"; NewText = \"4+94=\" },\n    @{ Row = 14; Col = 1; OldText = \"35-15=\"; NewText = \"3+58=\" },\n    @{ Row = 14; Col = 2; OldText = \"75-69=\"; NewText = \"10+55=\" },\n    @{ Row = 14; Col = 3; OldText = \"38-10=\"; NewText = \"77-64=\" },\n    @{ Row = 14; Col = 4; OldText = \"36+41=\"; NewText = \"94-20=\" },\n    @{ Row = 14; Col = 5; OldText = \"58+34=\"; NewText = \"45+30=\" },\n    @{ Row = 15; Col = 1; OldText = \"75+3=\"; NewText = \"16+7=\" },\n    @{ Row = 15; Col = 2; OldText = \"52+42=\"; NewText = \"99-38=\" },\n    @{ Row = 15; Col = 3; OldText = \"30+22=\"; NewText = \"96-41=\" },\n    @{ Row = 15; Col = 4; OldText = \"85-39=\"; NewText = \"40+13=\" },\n    @{ Row = 15; Col = 5; OldText = \"61-59=\"; NewText = \"55+13=\" },\n    @{ Row = 16; Col = 1; OldText = \"34+8=\"; NewText = \"32+51=\" },\n    @{ Row = 16; Col = 2; OldText = \"77+21=\"; NewText = \"73-59=\" },\n    @{ Row = 16; Col = 3; OldText = \"89-77=\"; NewText = \"88-30=\" },\n    @{ Row = 16; Col = 4; OldText = \"32+32=\"; NewText = \"26+13=\" },\n    @{ Row = 16; Col = 5; OldText = \"61+17=\"; NewText = \"51-18=\" },\n    @{ Row = 17; Col = 1; OldText = \"2+48=\"; NewText = \"36+40=\" },\n    @{ Row = 17; Col = 2; OldText = \"38-25=\"; NewText = \"11+1=\" },\n    @{ Row = 17; Col = 3; OldText = \"68-62=\"; NewText = \"36+62=\" },\n    @{ Row = 17; Col = 4; OldText = \"71-12=\"; NewText = \"23-18=\" },\n    @{ Row = 17; Col = 5; OldText = \"98-5=\"; NewText = \"98-93=\" },\n    @{ Row = 18; Col = 1; OldText = \"59-50=\"; NewText = \"14+19=\" },\n    @{ Row = 18; Col = 2; OldText = \"83-82=\"; NewText = \"11+8=\" },\n    @{ Row = 18; Col = 3; OldText = \"55+12=\"; NewText = \"10+3=\" },\n    @{ Row = 18; Col = 4; OldText = \"97-51=\"; NewText = \"96-52=\" },\n    @{ Row = 18; Col = 5; OldText = \"58+19=\"; NewText = \"20-13=\" },\n    @{ Row = 19; Col = 1; OldText = \"86-7=\"; NewText = \"55-0=\" },\n    @{ Row = 19; Col = 2; OldText = \"65+18=\"; NewText = \"22-15=\" },\n    @{ Row = 19; Col = 3; OldText = \"71-64=\"; NewText = \"51-29=\" },\n    @{ Row = 19; Col = 4; OldText = \"64-58=\"; NewText = \"31+30=\" },\n    @{ Row = 19; Col = 5; OldText = \"82-33=\"; NewText = \"63-4=\" },\n    @{ Row = 20; Col = 1; OldText = \"98-59=\"; NewText = \"23+52=\" },\n    @{ Row = 20; Col = 2; OldText = \"1+28=\"; NewText = \"28-7=\" },\n    @{ Row = 20; Col = 3; OldText = \"42+38=\"; NewText = \"15-14=\" },\n    @{ Row = 20; Col = 4; OldText = \"59+39=\"; NewText = \"42-39=\" },\n    @{ Row = 20; Col = 5; OldText = \"55-50=\"; NewText = \"92-85=\" },\n)\n\nforeach ($u in $cellUpdates) {\n    $cell = $table.Cell($u.Row, $u.Col)\n    $cellRange = $cell.Range\n    $currentText = $cellRange.Text.TrimEnd([char]13,[char]7)\n    if ($currentText -ne $u.OldText) {\n        throw \"Cell (\" + $u.Row + \",\" + $u.Col + \") expected '\" + $u.OldText + \"' but found '\" + $currentText + \"'\"\n    }\n    $cellRange.Text = $u.NewText\n}\n"}
